$wb = $excel.ActiveWorkbook

$rowValues = @{
    A = 12
    B = "2026-02-16"
    C = "22:53:10"
    D = "base_strategy"
    E = "DOWN"
    F = 49.999998
    G = $null
    H = "OPEN"
    I = 0
    J = 0
    K = 100
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = $null
    Q = 0
}

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A13").Value = $rowValues.A
    # Prefix with a quote so the date-shaped string is kept as literal text
    # instead of being auto-parsed into a date serial number; reset the
    # style afterwards so the quote-prefix formatting doesn't stick.
    $ws.Range("B13").Value = "'" + $rowValues.B
    $ws.Range("B13").Style = "Normal"
    $ws.Range("C13").Value = $rowValues.C
    $ws.Range("D13").Value = $rowValues.D
    $ws.Range("E13").Value = $rowValues.E
    $ws.Range("F13").Value = $rowValues.F
    # Quote-prefix trick also used for an explicit empty text cell (rather
    # than a cleared/blank cell) to match the other OPEN rows above it.
    $ws.Range("G13").Value = "'"
    $ws.Range("G13").Style = "Normal"
    $ws.Range("H13").Value = $rowValues.H
    $ws.Range("I13").Value = $rowValues.I
    $ws.Range("J13").Value = $rowValues.J
    $ws.Range("K13").Value = $rowValues.K
    $ws.Range("L13").Value = $rowValues.L
    $ws.Range("M13").Value = $rowValues.M
    $ws.Range("N13").Value = $rowValues.N
    $ws.Range("O13").Value = $rowValues.O
    $ws.Range("P13").Value = "'"
    $ws.Range("P13").Style = "Normal"
    $ws.Range("Q13").Value = $rowValues.Q
}
